# Fixed the scoring system:
# - Re-ordered / corrected the INPUT_SENTENCE, NAME and SENTENCES values per row
# - Updated the Calculated_Intelligibility (E) and Words_Correct (G) scores

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = "Enjoy the fair weather while in the tropics."
$ws.Range("C2").Value = "jewpofjoiwjFOJWROIFJWERFJOP[WKJEFPJKWF"
$ws.Range("D2").Value = "P1_W1_S4"
$ws.Range("E2").Value = 0.2439024390243902
$ws.Range("G2").Value = 8

# Row 3
$ws.Range("B3").Value = "You're used to being on the field."
$ws.Range("C3").Value = "PRKF[PJrwpvjwprjvp"
$ws.Range("D3").Value = "P1_W1_S3"
$ws.Range("E3").Value = 0.07692307692307693

# Row 4
$ws.Range("B4").Value = "We picked grapes for wine"
$ws.Range("C4").Value = "khiuh wfhouhf ojhfojq oejfojq oijfojw oijwefojw jwoejfoi oijowj feohweoh"
$ws.Range("D4").Value = "P1_W1_S1"
$ws.Range("E4").Value = 0.2061855670103093
$ws.Range("G4").Value = 9

# Row 5
$ws.Range("B5").Value = "The ballet is about to begin."
$ws.Range("C5").Value = "uhfowhoufh fweoij wefjowj ewfjojo efoijwo oijewoj efwjwo"
$ws.Range("D5").Value = "P1_W1_S2"
$ws.Range("E5").Value = 0.1882352941176471
$ws.Range("G5").Value = 7

# Row 6
$ws.Range("B6").Value = "he is capable and willing to make decisions."
$ws.Range("C6").Value = "owijefohj ohfweoh ewohfo efhweo fhjowhj fowejofj fjowej fewojo"
$ws.Range("D6").Value = "P1_W2_S4"
$ws.Range("E6").Value = 0.169811320754717
$ws.Range("G6").Value = 8

# Row 7
$ws.Range("B7").Value = "Big muscles are not necessarily strong ones"
$ws.Range("C7").Value = "whfi wefoh owjowj ojwfo fwjoj oefjoqj wfoijwo oewjoewj fjwojo wfjwojf ofewijfo"
$ws.Range("D7").Value = "P1_W2_S3"
$ws.Range("E7").Value = 0.115702479338843
$ws.Range("G7").Value = 11

# Row 8
$ws.Range("B8").Value = "I think I'm getting better."
$ws.Range("C8").Value = "fnk hfow wrfh weojfo fwojo wroijfho ojwro fojoiwr ojwof jwfoi wjojo"
$ws.Range("D8").Value = "P1_W2_S1"
$ws.Range("E8").Value = 0.148936170212766
$ws.Range("G8").Value = 11

# Row 9
$ws.Range("B9").Value = "You want him to do well"
$ws.Range("C9").Value = "ncaj ednfow woejfo jfowe weijfo pfjwoi foijwoi ewpijpo wpjfp ieoaij eijfp"
$ws.Range("D9").Value = "P1_W2_S2"
$ws.Range("E9").Value = 0.1875
$ws.Range("G9").Value = 11
